$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2"
$ws.Range("D2").Value = 0.08795
$ws.Range("E2").Value = 0.004300000000000012
$ws.Range("G2").Value = 0.0003724736930015033
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 127.8
$ws.Range("L2").Value = 0.213462502087857
$ws.Range("M2").Value = 41.63
$ws.Range("N2").Value = 0.02789466630930046
$ws.Range("O2").Value = 0.3257433489827856
$ws.Range("P2").Value = 40.2
$ws.Range("Q2").Value = 0.02693647815599035
$ws.Range("R2").Value = 0.3145539906103286
$ws.Range("S2").Value = 1.43
$ws.Range("T2").Value = 0.03435022820081671
$ws.Range("U2").Value = 532.4
$ws.Range("V2").Value = 0.3567408201554543
$ws.Range("W2").Value = 0.09405231931489444
$ws.Range("X2").Value = 0.02496830015975293
$ws.Range("Y2").Value = 0.06908401915514151
$ws.Range("Z2").Value = 0.2435026640094359
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 0.0209625367197816
$ws.Range("AC2").Value = -0.0209625367197816
$ws.Range("AD2").Value = 1568.7
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 1568.7
$ws.Range("AG2").Value = 1036.3
$ws.Range("AH2").Value = 0.5124628401554996
$ws.Range("AI2").Value = 0.5200053038087977
$ws.Range("AJ2").Value = 0.4098153201249654
$ws.Range("AK2").Value = 0.4171396369198567
$ws.Range("AN2").ClearContents()
$ws.Range("AP2").ClearContents()

# --- Row 3 updates ---
$ws.Range("D3").Value = 0.083
$ws.Range("E3").Value = 0.108
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 116.4
$ws.Range("L3").Value = 0.3708187320802804
$ws.Range("M3").Value = 41.63
$ws.Range("N3").Value = 0.03798704261337713
$ws.Range("O3").Value = 0.3576460481099656
$ws.Range("P3").Value = 40.2
$ws.Range("Q3").Value = 0.03668217903093348
$ws.Range("R3").Value = 0.345360824742268
$ws.Range("S3").Value = 1.43
$ws.Range("T3").Value = 0.03435022820081671
$ws.Range("U3").Value = 532.4
$ws.Range("V3").Value = 0.4858107491559448
$ws.Range("W3").Value = 0.1648725212464589
$ws.Range("X3").Value = 0.02173741790406666
$ws.Range("Y3").Value = 0.1431351033423923
$ws.Range("Z3").Value = 0.2973382589750876
$ws.Range("AA3").Value = 0
$ws.Range("AB3").Value = 0.01884756420890493
$ws.Range("AC3").Value = -0.01884756420890493
$ws.Range("AD3").Value = 843.1
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 843.1
$ws.Range("AG3").Value = 310.7
$ws.Range("AH3").Value = 0.4348117586384734
$ws.Range("AI3").Value = 0.4926376066378403
$ws.Range("AJ3").Value = 0.2208872458410351
$ws.Range("AK3").Value = 0.2635284139100934
$ws.Range("AN3").ClearContents()
$ws.Range("AP3").ClearContents()

# --- Row 4 updates (Avida -> Hoist Finance) ---
$ws.Range("B4").Value = "Hoist Finance AB (publ) (OM:HOFI)"
$ws.Range("D4").Value = 0.0929
$ws.Range("E4").Value = -0.09939999999999999
$ws.Range("G4").Value = 0.0007830056179775281
$ws.Range("K4").Value = 11.4
$ws.Range("L4").Value = 0.04002808988764045
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0.02323211738332994
$ws.Range("X4").Value = 0.0281991824154392
$ws.Range("Y4").Value = -0.004967065032109259
$ws.Range("Z4").Value = 0.2029935851746258
$ws.Range("AB4").Value = 0.02307750923065827
$ws.Range("AC4").Value = -0.02307750923065827
$ws.Range("AD4").Value = 725.6
$ws.Range("AF4").Value = 725.6
$ws.Range("AG4").Value = 725.6
$ws.Range("AH4").Value = 0.6466446840745033
$ws.Range("AI4").Value = 0.5558875354324676
$ws.Range("AJ4").Value = 0.6466446840745033
$ws.Range("AK4").Value = 0.5558875354324676

# --- Remove row 5 (merged into row 4) ---
$ws.Rows("5:5").Delete()
